$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.158.19'
$ws.Range("E2").Value = '  -0.15%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.842.12'
$ws.Range("E3").Value = '  -0.34%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9989'
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.39'
$ws.Range("E5").Value = '  -1.73%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6881'
$ws.Range("E6").Value = '  -1.98%  '
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3012'
$ws.Range("E8").Value = '  -1.91%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07456'
$ws.Range("E9").Value = '  -3.37%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.13'
$ws.Range("E10").Value = '  -2.33%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07654'
$ws.Range("E11").Value = '  -2.12%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.834.64'
$ws.Range("E12").Value = '  -0.85%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.065'
$ws.Range("E13").Value = '  -1.46%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6826'
$ws.Range("E14").Value = '  -0.52%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '87.52'
$ws.Range("E15").Value = '  -5.83%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.171'
$ws.Range("E16").Value = '  -6.59%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '29.144.89'
$ws.Range("E17").Value = '  -0.08%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008150'
$ws.Range("E18").Value = '  -1.89%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.078.21'
$ws.Range("E19").Value = '  -0.38%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '228.55'
$ws.Range("E20").Value = '  -5.65%  '
$ws.Range("E21").Value = '  -1.55%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9999'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.404'
$ws.Range("E23").Value = '  -1.57%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.9997'
$ws.Range("E24").Value = '  -0.07%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1455'
$ws.Range("E25").Value = '  -3.77%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '159.54'
$ws.Range("E26").Value = '  +0.15%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.778'
$ws.Range("E27").Value = '  -0.57%  '
$ws.Range("E28").Value = '  -1.27%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.513'
$ws.Range("E29").Value = '  -1.80%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.280'
$ws.Range("E30").Value = '  +1.24%  '
$ws.Range("E31").Value = '  -1.03%  '
$ws.Range("E32").Value = '  -0.95%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05260'
$ws.Range("E33").Value = '  +2.69%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7586'
$ws.Range("E34").Value = '  -4.26%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.853'
$ws.Range("E35").Value = '  -3.31%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.134'
$ws.Range("E36").Value = '  -1.14%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.687'
$ws.Range("E37").Value = '  -0.42%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.305.57'
$ws.Range("E38").Value = '  -1.35%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01834'
$ws.Range("E39").Value = '  -2.03%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.725'
$ws.Range("E40").Value = '  +0.37%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9305'
$ws.Range("E41").Value = '  -2.73%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.931'
$ws.Range("E42").Value = '  -2.24%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '104.98'
$ws.Range("E43").Value = '  -2.17%  '
$ws.Range("E44").Value = '  -0.12%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.980.31'
$ws.Range("E45").Value = '  -0.36%  '
$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5198'
$ws.Range("E46").Value = '  +0.29%  '
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '65.01'
$ws.Range("E47").Value = '  +0.98%  '
$ws.Range("E48").Value = '  -0.28%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.526'
$ws.Range("E49").Value = '  -1.94%  '
$ws.Range("E50").Value = '  +0.44%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05948'
$ws.Range("E51").Value = '  +0.74%  '
